# Fruta / hortaliza, semanal
# Feria Lagunitas de Puerto Montt - Nectarin: add a new weekly price record
# (Super Queen) on top of the sheet's last "June Pearl" pair, pushing that
# pair down two rows (the rest of the table below shifts down accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new blank rows at 359/360, shifting the former rows 359..372
#    down to 361..374 (dimension grows from A1:T372 to A1:T374).
$ws.Range("A359:A360").EntireRow.Insert()

# 2) The two new blank rows (359, 360) receive the data that used to live in
#    rows 357/358 (June Pearl, "Especial"/"Primera") before those rows are
#    overwritten with the new "Super Queen" entries below.
$ws.Range("A359").Value = 4
$ws.Range("B359").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C359").Value = "Los Lagos"
$ws.Range("D359").Value = 44266
$ws.Range("E359").Value = 10
$ws.Range("F359").Value = "Fruta"
$ws.Range("G359").Value = 100103
$ws.Range("H359").Value = "Frutos de hueso (carozo)"
$ws.Range("I359").Value = 100103006
$ws.Range("J359").Value = "Nectarín"
$ws.Range("K359").Value = "June Pearl"
$ws.Range("L359").Value = "Especial"
$ws.Range("M359").Value = 100
$ws.Range("N359").Value = 18000
$ws.Range("O359").Value = 18000
$ws.Range("P359").Value = 18000
$ws.Range("Q359").Value = "$/caja 15 kilos empedrada"
$ws.Range("R359").Value = "Región de O'Higgins"
$ws.Range("S359").Value = 1200
$ws.Range("T359").Value = 15

$ws.Range("A360").Value = 4
$ws.Range("B360").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C360").Value = "Los Lagos"
$ws.Range("D360").Value = 44266
$ws.Range("E360").Value = 10
$ws.Range("F360").Value = "Fruta"
$ws.Range("G360").Value = 100103
$ws.Range("H360").Value = "Frutos de hueso (carozo)"
$ws.Range("I360").Value = 100103006
$ws.Range("J360").Value = "Nectarín"
$ws.Range("K360").Value = "June Pearl"
$ws.Range("L360").Value = "Primera"
$ws.Range("M360").Value = 100
$ws.Range("N360").Value = 15000
$ws.Range("O360").Value = 15000
$ws.Range("P360").Value = 15000
$ws.Range("Q360").Value = "$/caja 15 kilos empedrada"
$ws.Range("R360").Value = "Región de O'Higgins"
$ws.Range("S360").Value = 1000
$ws.Range("T360").Value = 15

# 3) Overwrite rows 357/358 with the new "Super Queen" weekly records
#    (only D, K, M, N, O, P, S change; the rest of the row stays as-is).
$ws.Range("D357").Value = 44578
$ws.Range("K357").Value = "Super Queen"
$ws.Range("M357").Value = 200
$ws.Range("N357").Value = 20000
$ws.Range("O357").Value = 20000
$ws.Range("P357").Value = 20000
$ws.Range("S357").Value = 1333

$ws.Range("D358").Value = 44578
$ws.Range("K358").Value = "Super Queen"
$ws.Range("M358").Value = 400
$ws.Range("N358").Value = 16500
$ws.Range("O358").Value = 17000
$ws.Range("P358").Value = 16750
$ws.Range("S358").Value = 1117
